$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new to-do item as a new row in column A
$ws.Range("A11").Value = "maybe give every row a unique number in addition to 1-8 for all conditions"

# Update the selection/active cell to E10
$ws.Range("E10").Select()
